# Updated symbol list on Tue Feb  7 22:58:05 UTC 2023 with GitHub Actions
# Refresh Price (column D) and Volume(1h) (column E) figures for the cryptos sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "332.42"
Set-TextValue "E2" "1.57%"
Set-TextValue "D3" "44.72"
Set-TextValue "E3" "2.04%"
Set-TextValue "D4" "5.549"
Set-TextValue "E4" "-0.16%"
Set-TextValue "D5" "0.08283"
Set-TextValue "E5" "2.73%"
Set-TextValue "D6" "2.044"
Set-TextValue "E6" "3.72%"
Set-TextValue "D7" "0.9778"
Set-TextValue "E7" "3.35%"
Set-TextValue "D8" "0.1129"
Set-TextValue "E8" "-3.28%"
Set-TextValue "D9" "0.1913"
Set-TextValue "E9" "3.10%"
Set-TextValue "D10" "10.26"
Set-TextValue "E10" "-13.21%"
Set-TextValue "D11" "0.1008"
Set-TextValue "E11" "2.46%"
Set-TextValue "D12" "0.04678"
Set-TextValue "E12" "-1.84%"
Set-TextValue "D13" "0.1059"
Set-TextValue "E13" "-0.66%"
Set-TextValue "D14" "0.001264"
Set-TextValue "E14" "-1.82%"
Set-TextValue "E15" "-2.44%"
Set-TextValue "D16" "0.006028"
Set-TextValue "E16" "2.68%"
Set-TextValue "E17" "-0.29%"
Set-TextValue "D18" "4.435"
Set-TextValue "E18" "2.58%"
Set-TextValue "E20" "-3.47%"
Set-TextValue "D21" "0.1384"
Set-TextValue "E21" "-1.73%"
Set-TextValue "D22" "0.2489"
Set-TextValue "E22" "-0.77%"
Set-TextValue "D23" "0.001302"
Set-TextValue "E23" "4.24%"
Set-TextValue "D24" "0.004408"
Set-TextValue "E24" "2.36%"
Set-TextValue "D25" "0.0001279"
Set-TextValue "E25" "7.25%"
Set-TextValue "D26" "0.0003739"
Set-TextValue "E26" "-0.36%"
Set-TextValue "D38" "0.02809"
Set-TextValue "E38" "9.27%"
Set-TextValue "D39" "0.05763"
Set-TextValue "E39" "4.69%"
Set-TextValue "D40" "0.007643"
Set-TextValue "E40" "1.13%"
Set-TextValue "D41" "0.1428"
Set-TextValue "E41" "2.07%"
Set-TextValue "D42" "0.007545"
Set-TextValue "E42" "1.31%"
Set-TextValue "E43" "-2.35%"
Set-TextValue "D44" "0.008028"
Set-TextValue "E44" "-3.95%"
Set-TextValue "D45" "0.00007034"
Set-TextValue "E45" "-0.89%"
Set-TextValue "D46" "0.00000000750"
Set-TextValue "E46" "-0.24%"
Set-TextValue "D47" "0.0005798"
Set-TextValue "E47" "-0.23%"
Set-TextValue "D48" "0.003551"
Set-TextValue "E48" "-26.62%"
Set-TextValue "D50" "0.00002100"
Set-TextValue "E50" "-0.24%"
Set-TextValue "D51" "0.0002000"
Set-TextValue "E51" "-0.24%"
